$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '256.35'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.78%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '27.08'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.71%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.645'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-10.96%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05876'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.96%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.621'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-1.31%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8587'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.90%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9427'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-5.94%'
$ws.Range("B9").Value = 'One'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.01039'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1,617.54%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1404'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-1.10%'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.04380'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '22.87%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07102'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.07%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03150'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.23%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09144'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.88%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001523'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.49%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006226'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '5.66%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.520'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.64%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.205'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.94%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.58%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.821'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '8.43%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04242'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.38%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001220'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.14%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004286'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-5.28%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001201'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.10%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.05%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03821'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.39%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006214'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-5.53%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1100'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-0.58%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002201'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '0.10%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01141'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '5.15%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005466'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.62%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.10%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05003'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-54.12%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.2101'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '9,310.95%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.10%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.10%'
